$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A3").Value = "dadaaaaaaa"
$ws.Range("A3").Select()
